$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (Case_0_25, 380 kV slack voltage run)
$newValues = @{
    2 = @{ B=1.019999999999999; C=1.022897867616783; D=1.024548578315914; E=1.023566386591893; F=1.021353823986725; I=1.028836520122798; J=1.028081120431348; K=1.027376755783184; L=1.026397447707555; M=1.024191403148723; N=1.013345579854189 }
    3 = @{ B=1.02; C=1.024112492986348; D=1.025612430444894; E=1.024604717178594; F=1.023210306529709; I=1.02901050710983; J=1.028932617071243; K=1.028247024285459; L=1.02724205374203; M=1.02585144794885; N=1.0136334068263 }
    4 = @{ B=1.02; C=1.024897374572487; D=1.026300124278936; E=1.025276011632317; F=1.02441020130392; I=1.029121183798649; J=1.029482068952284; K=1.028808869025283; L=1.027787402178921; M=1.026923833129279; N=1.013818980799488 }
    5 = @{ B=1.02; C=1.025227088642129; D=1.026589068635823; E=1.025558088940975; F=1.024914320199654; I=1.029167256991269; J=1.029712696765717; K=1.02904476560129; L=1.028016389629528; M=1.027374249578921; N=1.013896836746027 }
    6 = @{ B=1.02; C=1.025282434536608; D=1.026637574243117; E=1.025605443132205; F=1.024998945745689; I=1.029174966203776; J=1.029751399061863; K=1.029084355990197; L=1.028054821460604; M=1.027449852496352; N=1.013909899783141 }
    7 = @{ B=1.02; C=1.024901781205811; D=1.026303985800989; E=1.025279781289733; F=1.024416938595591; I=1.029121801219071; J=1.029485152030261; K=1.02881202227383; L=1.027790463005771; M=1.026929853230686; N=1.013820021739514 }
    8 = @{ B=1.02; C=1.023308578401713; D=1.024908256710589; E=1.02391741596438; F=1.021981522584479; I=1.028895714584516; J=1.028369204912996; K=1.027671133241189; L=1.026683129729496; M=1.024752796952905; N=1.013442991450892 }
    9 = @{ B=1.02; C=1.020492828727199; D=1.022443387821461; E=1.021512227276311; F=1.017678950802846; I=1.028482704505776; J=1.026390963257555; K=1.025650832601068; L=1.024722801819336; M=1.020902475815382; N=1.01277344954745 }
    10 = @{ B=1.02; C=1.01860980188011; D=1.020796309888351; E=1.019905547881662; F=1.014802381645719; I=1.02819749344782; J=1.025064013057227; K=1.0242971236889; L=1.023409653109198; M=1.018325466697646; N=1.012323557447222 }
    11 = @{ B=1.02; C=1.017792980709856; D=1.020082154698731; E=1.019209033928194; F=1.013554685261694; I=1.028071642444004; J=1.02448745964746; K=1.023709289119862; L=1.022839521305497; M=1.017207041708533; N=1.012127898644732 }
    12 = @{ B=1.02; C=1.017489352641783; D=1.019816737659617; E=1.018950192196633; F=1.013090901784383; I=1.028024541359337; J=1.02427300149963; K=1.023490686717844; L=1.022627515745454; M=1.016791210876963; N=1.012055092997477 }
    13 = @{ B=1.02; C=1.017554492128643; D=1.01987367727502; E=1.019005720360669; F=1.013190400240271; I=1.028034660755397; J=1.024319017161608; K=1.02353758924893; L=1.022673002266473; M=1.016880426220592; N=1.012070715923997 }
    14 = @{ B=1.02; C=1.017767887325551; D=1.02006021828506; E=1.019187640569562; F=1.013516355662337; I=1.028067756289202; J=1.024469738623029; K=1.023691224589292; L=1.022822001637215; M=1.017172677196371; N=1.012121883150821 }
    15 = @{ B=1.02; C=1.017899337267353; D=1.02017513263052; E=1.019299710885415; F=1.013717142928573; I=1.02808810054746; J=1.024562563151414; K=1.023785850564479; L=1.022913774060633; M=1.017352689541593; N=1.012153391805066 }
    16 = @{ B=1.02; C=1.018663980443001; D=1.020843685472367; E=1.019951755824708; F=1.014885141179314; I=1.028205796114554; J=1.025102235070314; K=1.024336100846954; L=1.02344745836025; M=1.018399637779724; N=1.012336524615683 }
    17 = @{ B=1.02; C=1.019143226832241; D=1.021262791403321; E=1.020360546706955; F=1.015617217023849; I=1.028278992876208; J=1.025440225523704; K=1.024680808747352; L=1.023781812478016; M=1.019055665460718; N=1.012451169914746 }
    18 = @{ B=1.02; C=1.019422622813346; D=1.021507156479244; E=1.020598909489756; F=1.016044020631899; I=1.028321460369221; J=1.025637179315525; K=1.0248817101812; L=1.023976688211648; M=1.019438068890981; N=1.012517958398931 }
    19 = @{ B=1.02; C=1.019517865978145; D=1.021590463110778; E=1.020680171911847; F=1.016189515494597; I=1.028335902222269; J=1.02570430330068; K=1.024950185203658; L=1.024043110906776; M=1.019568417171184; N=1.012540717642127 }
    20 = @{ B=1.02; C=1.019091822784575; D=1.021217834881644; E=1.020316695411227; F=1.01553869340661; I=1.028271163037061; J=1.025403982061066; K=1.024643841516897; L=1.023745954713862; M=1.018985305483119; N=1.01243887806815 }
    21 = @{ B=1.02; C=1.017705053993418; D=1.020005290723669; E=1.019134073093792; F=1.013420379237975; I=1.028058020270727; J=1.024425363236254; K=1.023645989864674; L=1.022778131501784; M=1.017086627659061; N=1.012106819255397 }
    22 = @{ B=1.02; C=1.016831836955469; D=1.019242057851682; E=1.018389783971667; F=1.012086574306288; I=1.027921957920695; J=1.023808326003724; K=1.023017126824211; L=1.02216827157245; M=1.015890544409488; N=1.011897292274542 }
    23 = @{ B=1.02; C=1.017294870861402; D=1.019646744679743; E=1.018784415879549; F=1.012793838001866; I=1.027994281880271; J=1.024135595390434; K=1.023350639965662; L=1.022491699067323; M=1.016524833994; N=1.012008437827369 }
    24 = @{ B=1.02; C=1.019115050494305; D=1.021238149083707; E=1.020336510165099; F=1.015574175477667; I=1.028274701705503; J=1.025420359509561; K=1.024660545913884; L=1.023762157749135; M=1.019017098888523; N=1.012444432478384 }
    25 = @{ B=1.02; C=1.021221780912991; D=1.023081276832552; E=1.022134580735897; F=1.018792656344311; I=1.028591214099793; J=1.026903803836299; K=1.02617432135506; L=1.02523068508656; M=1.021899611988771; N=1.012947160271943 }
}

foreach ($r in $newValues.Keys) {
    foreach ($col in $newValues[$r].Keys) {
        $ws.Range("$col$r").Value = $newValues[$r][$col]
    }
}
